$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text instead of converting to a numeric value.
$riskyCells = @('D5','D6','D8','D13','D17','D19','D20','D21','D22','D23','D25','D26','D27','D28','D30','D32','D33','D35','D36','D37','D39','D40','D43','D44','D46','D49','D50','D51')
foreach ($addr in $riskyCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + volume change), including the
# reordering of BitcoinCash/SuiNetwork and VeChain/MantraDAO rows.
$ws.Range('D2').Value = '98.752.05'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '3.350.75'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '260.12'
$ws.Range('D6').Value = '653.82'
$ws.Range('E6').Value = '  +3.89%  '
$ws.Range('E7').Value = '  +12.77%  '
$ws.Range('D8').Value = '0.464'
$ws.Range('E8').Value = '  +18.02%  '
$ws.Range('E9').Value = '  +26.16%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.348.14'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D13').Value = '43.59'
$ws.Range('E13').Value = '  +21.05%  '
$ws.Range('E14').Value = '  +8.75%  '
$ws.Range('D15').Value = '98.584.06'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').Value = '3.985.22'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '5.60'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '3.353.23'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '7.56'
$ws.Range('E19').Value = '  +23.10%  '
$ws.Range('D20').Value = '16.93'
$ws.Range('E20').Value = '  +10.55%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '535.95'
$ws.Range('E21').Value = '  +8.57%  '
$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').Value = '3.61'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').Value = '10.41'
$ws.Range('E23').Value = '  +9.89%  '
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').Value = '0.437'
$ws.Range('E25').Value = '  +57.80%  '
$ws.Range('D26').Value = '101.82'
$ws.Range('E26').Value = '  +14.80%  '
$ws.Range('D27').Value = '6.27'
$ws.Range('E27').Value = '  +10.54%  '
$ws.Range('D28').Value = '12.63'
$ws.Range('E28').Value = '  +5.68%  '
$ws.Range('D29').Value = '3.528.19'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '0.150'
$ws.Range('E30').Value = '  +14.66%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '11.06'
$ws.Range('E32').Value = '  +15.66%  '
$ws.Range('D33').Value = '0.192'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').Value = '29.43'
$ws.Range('E35').Value = '  +5.69%  '
$ws.Range('D36').Value = '0.539'
$ws.Range('E36').Value = '  +16.90%  '
$ws.Range('D37').Value = '7.92'
$ws.Range('E37').Value = '  +7.84%  '
$ws.Range('E38').Value = '  +6.75%  '
$ws.Range('D39').Value = '0.157'
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('D40').Value = '524.17'
$ws.Range('E40').Value = '  +4.86%  '
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0432'
$ws.Range('E43').Value = '  +32.42%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').Value = '3.77'
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('E45').Value = '  +3.02%  '
$ws.Range('D46').Value = '0.827'
$ws.Range('E46').Value = '  +5.41%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  +20.18%  '
$ws.Range('D49').Value = '2.06'
$ws.Range('E49').Value = '  +5.80%  '
$ws.Range('D50').Value = '5.12'
$ws.Range('E50').Value = '  +9.97%  '
$ws.Range('D51').Value = '164.57'
$ws.Range('E51').Value = '  +2.57%  '
